$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: Merge the "CREADOR: " run and the "ALEXANDER BARRIOS" run
# into a single run "CREADOR: ALEXANDER BARRIOS", leaving the
# _GoBack bookmark at the end of the paragraph.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(2)
$r = $p.Range

# Delete the text of the second run ("ALEXANDER BARRIOS"), which sits
# after the bookmark.
$delRange = $d.Range($r.Start, $r.End)
$delRange.Find.Execute("ALEXANDER BARRIOS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delRange.Text = ""

# Replace the first run's text so it reads "CREADOR: ALEXANDER BARRIOS".
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$findRange = $d.Range($r2.Start, $r2.End)
$findRange.Find.Execute("CREADOR: ", $true, $false, $false, $false, $false, $true, 1, $false, "CREADOR: ALEXANDER BARRIOS", 2) | Out-Null

# ------------------------------------------------------------------
# Part 2: Insert a new paragraph right after the "CREADOR" paragraph
# containing "Prueba práctica de reigh" (bold, en-US), with proofErr
# spell-check markers around each unrecognized word.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(2)
$p3.Range.InsertParagraphAfter() | Out-Null

$newPar = $d.Paragraphs.Item(3)
$newRange = $newPar.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Prueba</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>práctica</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> de </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>reigh</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>' + `
  '<w:sectPr w:rsidR="00000000"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/></w:sectPr>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml) | Out-Null
